$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 2.204947
$ws.Cells.Item(2, 8).Value = 6.614841
$ws.Cells.Item(2, 9).Value = 0.03384377946268709
$ws.Cells.Item(2, 10).Value = 0.03400382310856976
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.1375686666666667
$ws.Cells.Item(2, 14).Value = 0.412706
$ws.Cells.Item(2, 15).Value = 0.2896572731203081
$ws.Cells.Item(2, 16).Value = 0.2896572731203081
$ws.Cells.Item(2, 17).Value = 0.3033316188606667
$ws.Cells.Item(2, 18).Value = 2.729984569746
$ws.Cells.Item(2, 19).Value = 0.00980309687124703
$ws.Cells.Item(2, 20).Value = 0.009849454677293637

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 2.204947
$ws.Cells.Item(3, 8).Value = 6.614841
$ws.Cells.Item(3, 9).Value = 0.03384377946268709
$ws.Cells.Item(3, 10).Value = 0.03400382310856976
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.3373673333333334
$ws.Cells.Item(3, 14).Value = 1.012102
$ws.Cells.Item(3, 15).Value = 0.7103427268796919
$ws.Cells.Item(3, 16).Value = 0.7103427268796919
$ws.Cells.Item(3, 17).Value = 0.7438770895313335
$ws.Cells.Item(3, 18).Value = 6.694893805782001
$ws.Cells.Item(3, 19).Value = 0.02404068259144006
$ws.Cells.Item(3, 20).Value = 0.02415436843127612

# Row 4
$ws.Cells.Item(4, 5).Value = 3.0
$ws.Cells.Item(4, 6).Value = 1.0
$ws.Cells.Item(4, 7).Value = 7.466229000000001
$ws.Cells.Item(4, 8).Value = 22.398687
$ws.Cells.Item(4, 9).Value = 0.1145993113185572
$ws.Cells.Item(4, 10).Value = 0.1151412393150827
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.1375686666666667
$ws.Cells.Item(4, 14).Value = 0.412706
$ws.Cells.Item(4, 15).Value = 0.2896572731203081
$ws.Cells.Item(4, 16).Value = 0.2896572731203081
$ws.Cells.Item(4, 17).Value = 1.027119168558
$ws.Cells.Item(4, 18).Value = 9.244072517022001
$ws.Cells.Item(4, 19).Value = 0.03319452401799854
$ws.Cells.Item(4, 20).Value = 0.03335149740369968

# Row 5
$ws.Cells.Item(5, 5).Value = 3.0
$ws.Cells.Item(5, 6).Value = 1.0
$ws.Cells.Item(5, 7).Value = 7.466229000000001
$ws.Cells.Item(5, 8).Value = 22.398687
$ws.Cells.Item(5, 9).Value = 0.1145993113185572
$ws.Cells.Item(5, 10).Value = 0.1151412393150827
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.3373673333333334
$ws.Cells.Item(5, 14).Value = 1.012102
$ws.Cells.Item(5, 15).Value = 0.7103427268796919
$ws.Cells.Item(5, 16).Value = 0.7103427268796919
$ws.Cells.Item(5, 17).Value = 2.518861767786
$ws.Cells.Item(5, 18).Value = 22.669755910074
$ws.Cells.Item(5, 19).Value = 0.08140478730055867
$ws.Cells.Item(5, 20).Value = 0.08178974191138305

# Row 6
$ws.Cells.Item(6, 5).Value = 3.0
$ws.Cells.Item(6, 6).Value = 1.0
$ws.Cells.Item(6, 7).Value = 31.62082666666666
$ws.Cells.Item(6, 8).Value = 94.86247999999999
$ws.Cells.Item(6, 9).Value = 0.4853487562896166
$ws.Cells.Item(6, 10).Value = 0.4876439191146448
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.1375686666666667
$ws.Cells.Item(6, 14).Value = 0.412706
$ws.Cells.Item(6, 15).Value = 0.2896572731203081
$ws.Cells.Item(6, 16).Value = 0.2896572731203081
$ws.Cells.Item(6, 17).Value = 4.350034963431111
$ws.Cells.Item(6, 18).Value = 39.15031467088
$ws.Cells.Item(6, 19).Value = 0.1405847972591834
$ws.Cells.Item(6, 20).Value = 0.1412496078644481

# Row 7
$ws.Cells.Item(7, 5).Value = 3.0
$ws.Cells.Item(7, 6).Value = 1.0
$ws.Cells.Item(7, 7).Value = 31.62082666666666
$ws.Cells.Item(7, 8).Value = 94.86247999999999
$ws.Cells.Item(7, 9).Value = 0.4853487562896166
$ws.Cells.Item(7, 10).Value = 0.4876439191146448
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.3373673333333334
$ws.Cells.Item(7, 14).Value = 1.012102
$ws.Cells.Item(7, 15).Value = 0.7103427268796919
$ws.Cells.Item(7, 16).Value = 0.7103427268796919
$ws.Cells.Item(7, 17).Value = 10.66783397032889
$ws.Cells.Item(7, 18).Value = 96.01050573296
$ws.Cells.Item(7, 19).Value = 0.3447639590304333
$ws.Cells.Item(7, 20).Value = 0.3463943112501967

# Row 8
$ws.Cells.Item(8, 5).Value = 2.0
$ws.Cells.Item(8, 6).Value = 1.0
$ws.Cells.Item(8, 7).Value = 0.9199225
$ws.Cells.Item(8, 8).Value = 1.839845
$ws.Cells.Item(8, 9).Value = 0.01411991046168627
$ws.Cells.Item(8, 10).Value = 0.009457788014433987
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.1375686666666667
$ws.Cells.Item(8, 14).Value = 0.412706
$ws.Cells.Item(8, 15).Value = 0.2896572731203081
$ws.Cells.Item(8, 16).Value = 0.2896572731203081
$ws.Cells.Item(8, 17).Value = 0.1265525117616667
$ws.Cells.Item(8, 18).Value = 0.75931507057
$ws.Cells.Item(8, 19).Value = 0.004089934761034956
$ws.Cells.Item(8, 20).Value = 0.002739517086010882

# Row 9
$ws.Cells.Item(9, 5).Value = 2.0
$ws.Cells.Item(9, 6).Value = 1.0
$ws.Cells.Item(9, 7).Value = 0.9199225
$ws.Cells.Item(9, 8).Value = 1.839845
$ws.Cells.Item(9, 9).Value = 0.01411991046168627
$ws.Cells.Item(9, 10).Value = 0.009457788014433987
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.3373673333333334
$ws.Cells.Item(9, 14).Value = 1.012102
$ws.Cells.Item(9, 15).Value = 0.7103427268796919
$ws.Cells.Item(9, 16).Value = 0.7103427268796919
$ws.Cells.Item(9, 17).Value = 0.3103518006983333
$ws.Cells.Item(9, 18).Value = 1.86211080419
$ws.Cells.Item(9, 19).Value = 0.01002997570065131
$ws.Cells.Item(9, 20).Value = 0.006718270928423105

# Row 10
$ws.Cells.Item(10, 5).Value = 3.0
$ws.Cells.Item(10, 6).Value = 1.0
$ws.Cells.Item(10, 7).Value = 22.93880666666666
$ws.Cells.Item(10, 8).Value = 68.81642
$ws.Cells.Item(10, 9).Value = 0.3520882424674529
$ws.Cells.Item(10, 10).Value = 0.3537532304472688
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.1375686666666667
$ws.Cells.Item(10, 14).Value = 0.412706
$ws.Cells.Item(10, 15).Value = 0.2896572731203081
$ws.Cells.Item(10, 16).Value = 0.2896572731203081
$ws.Cells.Item(10, 17).Value = 3.155661048057778
$ws.Cells.Item(10, 18).Value = 28.40094943252
$ws.Cells.Item(10, 19).Value = 0.1019849202108443
$ws.Cells.Item(10, 20).Value = 0.1024671960888558

# Row 11
$ws.Cells.Item(11, 5).Value = 3.0
$ws.Cells.Item(11, 6).Value = 1.0
$ws.Cells.Item(11, 7).Value = 22.93880666666666
$ws.Cells.Item(11, 8).Value = 68.81642
$ws.Cells.Item(11, 9).Value = 0.3520882424674529
$ws.Cells.Item(11, 10).Value = 0.3537532304472688
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.3373673333333334
$ws.Cells.Item(11, 14).Value = 1.012102
$ws.Cells.Item(11, 15).Value = 0.7103427268796919
$ws.Cells.Item(11, 16).Value = 0.7103427268796919
$ws.Cells.Item(11, 17).Value = 7.738804034982222
$ws.Cells.Item(11, 18).Value = 69.64923631484
$ws.Cells.Item(11, 19).Value = 0.2501033222566086
$ws.Cells.Item(11, 20).Value = 0.251286034358413
